# feat - add scenario where correlated features are removed before feature selection
#
# Adds three new header columns (D:F) mirroring A:C ("rfe", "mutual
# information", "random forest") representing the new "correlated features
# removed" scenario, fills in the corresponding feature-selection results for
# D2:F11, updates the A:C values for the already-existing scenario, and
# clears out the rows/cells that are no longer used (rows 12-16 and a few
# trailing cells in columns A-C) so the sheet's used range shrinks from
# A1:C16 down to A1:F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (D1:F1), copying the bold/boxed style from A1:C1 ---
$ws.Range("A1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)

$ws.Range("D1").Value = "rfe"
$ws.Range("E1").Value = "mutual information"
$ws.Range("F1").Value = "random forest"

# --- Column A updates ---
$ws.Range("A3").Value = "Cell-Nominal-Voltage-V"
$ws.Range("A4").Value = "Cell-Energy-Wh"
$ws.Range("A5").Value = "Trigger-Mechanism_Nail"
$ws.Range("A6").Value = "Chemistry_NCA/Graphite"
$ws.Range("A7").Value = "Pressure-Assisted-Seal-Configuration-Negative"
$ws.Range("A8").Value = "Bottom-Vent-Yes-No"
$ws.Range("A9").ClearContents()
$ws.Range("A10").ClearContents()
$ws.Range("A11").ClearContents()
$ws.Range("A12").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("A14").ClearContents()
$ws.Range("A15").ClearContents()
$ws.Range("A16").ClearContents()

# --- Column B updates ---
$ws.Range("B2").Value = "Cell-Capacity-Ah"
$ws.Range("B3").Value = "Cell-Energy-Wh"
$ws.Range("B6").Value = "Bottom-Vent-Yes-No"
$ws.Range("B7").Value = "Chemistry_NMC/Graphite"
$ws.Range("B8").Value = "Trigger-Mechanism_Nail"
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("B13").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("B15").ClearContents()
$ws.Range("B16").ClearContents()

# --- Column C updates ---
$ws.Range("C2").Value = "Pre-Test-Cell-Mass-g"
$ws.Range("C3").Value = "Bottom-Vent-Yes-No"
$ws.Range("C4").Value = "Cell-Energy-Wh"
$ws.Range("C5").Value = "Trigger-Mechanism_Nail"
$ws.Range("C6").Value = "Cell-Capacity-Ah"
$ws.Range("C7").Value = "Cell-Nominal-Voltage-V"
$ws.Range("C8").Value = "Trigger-Mechanism_Heater (Non-ISC)"
$ws.Range("C9").Value = "Pressure-Assisted-Seal-Configuration-Positive"
$ws.Range("C10").Value = "Pressure-Assisted-Seal-Configuration-Negative"
$ws.Range("C11").Value = "Chemistry_NCA/Graphite"
$ws.Range("C12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("C16").ClearContents()

# --- New column D (rfe, correlated features removed) ---
$ws.Range("D2").Value = "Cell-Capacity-Ah"
$ws.Range("D3").Value = "Pre-Test-Cell-Mass-g"
$ws.Range("D4").Value = "Trigger-Mechanism_Nail"
$ws.Range("D5").Value = "Chemistry_NMC/Graphite"
$ws.Range("D6").Value = "Pressure-Assisted-Seal-Configuration-Negative"
$ws.Range("D7").Value = "Bottom-Vent-Yes-No"

# --- New column E (mutual information, correlated features removed) ---
$ws.Range("E2").Value = "Cell-Capacity-Ah"
$ws.Range("E3").Value = "Pre-Test-Cell-Mass-g"
$ws.Range("E4").Value = "Cell-Nominal-Voltage-V"
$ws.Range("E5").Value = "Bottom-Vent-Yes-No"
$ws.Range("E6").Value = "Chemistry_NMC/Graphite"
$ws.Range("E7").Value = "Trigger-Mechanism_Nail"

# --- New column F (random forest, correlated features removed) ---
$ws.Range("F2").Value = "Pre-Test-Cell-Mass-g"
$ws.Range("F3").Value = "Bottom-Vent-Yes-No"
$ws.Range("F4").Value = "Cell-Capacity-Ah"
$ws.Range("F5").Value = "Trigger-Mechanism_Nail"
$ws.Range("F6").Value = "Cell-Nominal-Voltage-V"
$ws.Range("F7").Value = "Trigger-Mechanism_Heater (Non-ISC)"
$ws.Range("F8").Value = "Pressure-Assisted-Seal-Configuration-Positive"
$ws.Range("F9").Value = "Pressure-Assisted-Seal-Configuration-Negative"
